$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the six pairing-time values in the second matrix (rows 14-20)
# from 240 minutes to 330 minutes. The matrix is symmetric, so each pair
# of collaborators has their shared total updated in both cells.
$ws.Range("F14").Value = 330   # Felix Schmidt <-> Magdalena Hinterkoerner
$ws.Range("E16").Value = 330   # Michael Baier <-> Jakob Stanta
$ws.Range("D17").Value = 330   # Jakob Stanta <-> Michael Baier
$ws.Range("B18").Value = 330   # Magdalena Hinterkoerner <-> Felix Schmidt
$ws.Range("H19").Value = 330   # Thomas Pinheiro de Souza <-> Florian Buchacher
$ws.Range("G20").Value = 330   # Florian Buchacher <-> Thomas Pinheiro de Souza

# Move the active selection to J24 to match the saved cursor position.
$ws.Range("J24").Select()
